# Update date format in Excel export to dd/mm/yy
# Column D currently holds date strings like "May 26, 2025" as text.
# Replace them with Excel serial date numbers and apply a dd/mm/yy number format.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row number -> Excel serial date value (matching the original text dates)
$dateValues = @{
    2  = 45803
    3  = 45804
    4  = 45805
    5  = 45806
    6  = 45807
    7  = 45810
    8  = 45811
    9  = 45812
    10 = 45813
    11 = 45814
    12 = 45817
    13 = 45818
    14 = 45819
    15 = 45820
    16 = 45821
    17 = 45824
    18 = 45825
    19 = 45826
    20 = 45827
    21 = 45828
}

foreach ($row in 2..21) {
    $cell = $ws.Cells.Item($row, 4)
    $cell.Value = $dateValues[$row]
    $cell.NumberFormat = "dd/mm/yy"
}

# Move the active selection from the old Task column (B2:B21) to the Date column (D2:D21)
$ws.Range("D2:D21").Select()
